$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 6897.2144
$ws.Range("J69").Value = 8815
$ws.Range("L69").Value = 26445
$ws.Range("N69").Value = -28193
$ws.Range("H72").Value = 6897.2144
$ws.Range("J72").Value = 8815
$ws.Range("L72").Value = 79335
$ws.Range("N72").Value = -88071
$ws.Range("H97").Value = 489.66666
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H118").Value = 1099.8182
$ws.Range("J118").Value = 1879
$ws.Range("L118").Value = 5637
$ws.Range("N118").Value = -8951
$ws.Range("H135").Value = 1148.6364
$ws.Range("I135").Value = 940.5
$ws.Range("J135").Value = 1703.6666
$ws.Range("K135").Value = 8464.5
$ws.Range("L135").Value = 15332.9994
$ws.Range("M135").Value = -5929.5
$ws.Range("N135").Value = -20402.9994
$ws.Range("H138").Value = 3301.842
$ws.Range("I138").Value = 5396.4
$ws.Range("J138").Value = 2984.4849
$ws.Range("K138").Value = 16189.2
$ws.Range("L138").Value = 8953.4547
$ws.Range("M138").Value = -11049.2
$ws.Range("N138").Value = -19233.4547

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4901.574
$ws.Range("I61").Value = 4894.28
$ws.Range("J61").Value = 4992.75
$ws.Range("K61").Value = 4894.28
$ws.Range("L61").Value = 4992.75
$ws.Range("M61").Value = -4682.28
$ws.Range("N61").Value = -5416.75
$ws.Range("H96").Value = 45308.6
$ws.Range("J96").Value = 45308.6
$ws.Range("L96").Value = 45308.6
$ws.Range("N96").Value = -50800.6
$ws.Range("H102").Value = 2544.375
$ws.Range("I102").Value = 2544.375
$ws.Range("K102").Value = 2544.375
$ws.Range("M102").Value = -922.375
$ws.Range("H136").Value = 4901.574
$ws.Range("I136").Value = 4894.28
$ws.Range("J136").Value = 4992.75
$ws.Range("K136").Value = 14682.84
$ws.Range("L136").Value = 14978.25
$ws.Range("M136").Value = -12132.84
$ws.Range("N136").Value = -20078.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2218.2
$ws.Range("I99").Value = 2122.8572
$ws.Range("K99").Value = 2122.8572
$ws.Range("M99").Value = -624.8571999999999
$ws.Range("H134").Value = 2049.3
$ws.Range("I134").Value = 1366.9474
$ws.Range("K134").Value = 4100.8422
$ws.Range("M134").Value = -1565.8422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 121671.11
$ws.Range("I31").Value = 6506
$ws.Range("K31").Value = 6506
$ws.Range("M31").Value = -6211
$ws.Range("H34").Value = 121671.11
$ws.Range("I34").Value = 6506
$ws.Range("K34").Value = 6506
$ws.Range("M34").Value = -6304
$ws.Range("H99").Value = 3800
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 3800
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 3800
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -6796
$ws.Range("H126").Value = 3800
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 11400
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -16340
$ws.Range("H131").Value = 25250
$ws.Range("J131").Value = 25000
$ws.Range("L131").Value = 25000
$ws.Range("N131").Value = -35080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 5750
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 15000
$ws.Range("M54").Value = -14441

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 13101310
$ws.Range("I132").Value = 13101310
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 39303930
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -39301400
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6456.0884
$ws.Range("I7").Value = 3864.24
$ws.Range("J7").Value = 13655.667
$ws.Range("K7").Value = 3864.24
$ws.Range("L7").Value = 13655.667
$ws.Range("M7").Value = -3752.24
$ws.Range("N7").Value = -13879.667
$ws.Range("H22").Value = 9391.471
$ws.Range("J22").Value = 10358
$ws.Range("L22").Value = 10358
$ws.Range("N22").Value = -10948
$ws.Range("H27").Value = 9391.471
$ws.Range("J27").Value = 10358
$ws.Range("L27").Value = 10358
$ws.Range("N27").Value = -10572
$ws.Range("H93").Value = 2098.5625
$ws.Range("I93").Value = 1948.6
$ws.Range("J93").Value = 2348.5
$ws.Range("K93").Value = 1948.6
$ws.Range("L93").Value = 2348.5
$ws.Range("M93").Value = -700.5999999999999
$ws.Range("N93").Value = -4844.5
$ws.Range("H122").Value = 8028.3
$ws.Range("I122").Value = 4571.5
$ws.Range("J122").Value = 10332.833
$ws.Range("K122").Value = 13714.5
$ws.Range("L122").Value = 30998.499
$ws.Range("M122").Value = -11264.5
$ws.Range("N122").Value = -35898.499
$ws.Range("H126").Value = 6456.0884
$ws.Range("I126").Value = 3864.24
$ws.Range("J126").Value = 13655.667
$ws.Range("K126").Value = 11592.72
$ws.Range("L126").Value = 40967.001
$ws.Range("M126").Value = -9122.719999999999
$ws.Range("N126").Value = -45907.001
$ws.Range("H136").Value = 5626.263
$ws.Range("I136").Value = 2145.9614
$ws.Range("K136").Value = 6437.8842
$ws.Range("M136").Value = -3887.8842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13971.143
$ws.Range("I41").Value = 10500
$ws.Range("J41").Value = 15359.6
$ws.Range("K41").Value = 10500
$ws.Range("L41").Value = 15359.6
$ws.Range("M41").Value = -10110
$ws.Range("N41").Value = -16139.6
$ws.Range("H62").Value = 7349.75
$ws.Range("J62").Value = 7500
$ws.Range("L62").Value = 7500
$ws.Range("N62").Value = -8748
$ws.Range("H65").Value = 7349.75
$ws.Range("J65").Value = 7500
$ws.Range("L65").Value = 37500
$ws.Range("N65").Value = -43740
$ws.Range("H94").Value = 64490
$ws.Range("J94").Value = 64490
$ws.Range("L94").Value = 64490
$ws.Range("N94").Value = -66292
$ws.Range("H122").Value = 4002.5625
$ws.Range("J122").Value = 8400.799999999999
$ws.Range("L122").Value = 25202.4
$ws.Range("N122").Value = -30102.4
